# Update cryptos list prices/volumes and fix row order for Aave / BitcoinSV
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) is stored as plain text in the source data (e.g.
# "69.20", "1.621.30"). Each Price cell being updated gets its
# NumberFormat forced to Text immediately before the value is written,
# so Excel's smart-typing doesn't reinterpret numeric-looking values as
# Numbers and silently drop significant trailing zeros (e.g.
# "69.20" -> 69.2). The Volume(1h) column (E) always contains a "%"
# sign plus surrounding spaces, so it is never at risk of numeric
# reinterpretation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.891.10"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.623.75"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.993"
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.62"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.29"
$ws.Range("E8").Value = "  +9.32%  "
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.856.81"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.621.30"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("E14").Value = "  +6.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.92"
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.920.13"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.82"
$ws.Range("E17").Value = "  +15.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.34"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.88"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("E22").Value = "  +2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.59"
$ws.Range("E23").Value = "  +4.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.45"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.61"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.58"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.423.06"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.65"
$ws.Range("E35").Value = "  +7.10%  "
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.555"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0505"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("E43").Value = "  +3.90%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "53.77"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.20"
$ws.Range("E45").Value = "  +5.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("E46").Value = "  +19.13%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.765.20"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.48"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0110"
$ws.Range("E51").Value = "  +9.90%  "
